# Add a new "Save" column (H) to the s_vals sheet, mirroring the header
# formatting used by the existing columns (e.g. "sum" in G1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last header cell (G1) onto the new header
# cell (H1) so it picks up the same bold/border/alignment style, then set
# its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" values for each data row.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
